# Updates the "cryptos" price list on Sheet1: refreshes the Price (column D)
# and Volume(1h) (column E) values for rows 2-51 to match the latest
# GitHub Actions data pull.
#
# Note: several Price values are plain numeric-looking strings (e.g. "1.001",
# "0.9998", "7.602"). Assigning such strings straight to Range.Value lets
# Excel auto-convert them to real numbers (losing formatting such as
# trailing zeros, e.g. "7.590" -> 7.59). To keep them as text exactly as in
# the source data, the cell's NumberFormat is forced to Text ("@") before
# the value is written for those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.120.34'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '1.748.60'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.56'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5551'
$ws.Range("E6").Value = '  +7.03%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2829'
$ws.Range("E8").Value = '  +0.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06183'
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("D10").Value = '1.753.46'
$ws.Range("E10").Value = '  +0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07224'
$ws.Range("E11").Value = '  +3.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.51'
$ws.Range("E12").Value = '  +0.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6533'
$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.644'
$ws.Range("E14").Value = '  +2.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.36'
$ws.Range("E15").Value = '  +1.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9999'
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").Value = '26.021.07'
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.81'
$ws.Range("E19").Value = '  +2.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006796'
$ws.Range("E20").Value = '  +2.82%  '

$ws.Range("D21").Value = '1.979.22'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("E22").Value = '  +5.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.760'
$ws.Range("E23").Value = '  +1.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.281'
$ws.Range("E24").Value = '  +2.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.97'
$ws.Range("E25").Value = '  +0.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.527'
$ws.Range("E26").Value = '  +1.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.33'
$ws.Range("E27").Value = '  +1.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.805'
$ws.Range("E28").Value = '  -1.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.35'
$ws.Range("E29").Value = '  +2.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08439'
$ws.Range("E30").Value = '  +2.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.799'
$ws.Range("E31").Value = '  +3.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.636'
$ws.Range("E32").Value = '  +5.88%  '

$ws.Range("E33").Value = '  +3.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.650'
$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.007'
$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6297'
$ws.Range("E36").Value = '  +2.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.706'
$ws.Range("E37").Value = '  +1.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01618'
$ws.Range("E38").Value = '  +1.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.972'
$ws.Range("E39").Value = '  +1.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9998'
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.42'
$ws.Range("E41").Value = '  +0.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3942'
$ws.Range("E42").Value = '  +2.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7477'
$ws.Range("E43").Value = '  +0.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.079'
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1153'
$ws.Range("E45").Value = '  +2.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.358'
$ws.Range("E46").Value = '  +0.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05333'
$ws.Range("E47").Value = '  -2.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.81'
$ws.Range("E48").Value = '  +3.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.89'
$ws.Range("E49").Value = '  +2.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3495'
$ws.Range("E50").Value = '  +1.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.602'
$ws.Range("E51").Value = '  -0.45%  '
